$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 163, shifting existing rows 163:182 down to 165:184
$ws.Rows.Item(163).Resize(2).Insert()

# New row 163 - Segunda quality entry
$ws.Range("A163").Value = 5
$ws.Range("B163").Value = "Macroferia Regional de Talca"
$ws.Range("C163").Value = "Maule"
$ws.Range("D163").Value = 45223
$ws.Range("E163").Value = 7
$ws.Range("F163").Value = "Fruta"
$ws.Range("G163").Value = 100107
$ws.Range("H163").Value = "Otros"
$ws.Range("I163").Value = 100107002
$ws.Range("J163").Value = "Chirimoya"
$ws.Range("K163").Value = "Cultivar IV Región"
$ws.Range("L163").Value = "Segunda"
$ws.Range("M163").Value = 230
$ws.Range("N163").Value = 19000
$ws.Range("O163").Value = 19000
$ws.Range("P163").Value = 19000
$ws.Range("Q163").Value = "`$/bandeja 10 kilos"
$ws.Range("R163").Value = "Provincia de Limarí"
$ws.Range("S163").Value = 1900
$ws.Range("T163").Value = 10

# New row 164 - Tercera quality entry
$ws.Range("A164").Value = 5
$ws.Range("B164").Value = "Macroferia Regional de Talca"
$ws.Range("C164").Value = "Maule"
$ws.Range("D164").Value = 45223
$ws.Range("E164").Value = 7
$ws.Range("F164").Value = "Fruta"
$ws.Range("G164").Value = 100107
$ws.Range("H164").Value = "Otros"
$ws.Range("I164").Value = 100107002
$ws.Range("J164").Value = "Chirimoya"
$ws.Range("K164").Value = "Cultivar IV Región"
$ws.Range("L164").Value = "Tercera"
$ws.Range("M164").Value = 180
$ws.Range("N164").Value = 17000
$ws.Range("O164").Value = 17000
$ws.Range("P164").Value = 17000
$ws.Range("Q164").Value = "`$/bandeja 10 kilos"
$ws.Range("R164").Value = "Provincia de Limarí"
$ws.Range("S164").Value = 1700
$ws.Range("T164").Value = 10
